# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet,
# mirroring the formatting of the existing "junio" column (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy the formatting (style) of column F into the new column G.
# Rows 1-6 share style s="2"/s="1", row 7 (totals row) uses style s="6".
$ws.Range("F1:F6").Copy()
$ws.Range("G1:G6").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the width of the new column to match the target (stored width 17).
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668

# Populate the new column's values.
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
